$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 17; existing rows 17-20 shift down to 18-21.
$ws.Rows.Item(17).Insert()

# Populate the new row 17 with a fresh weekly record (copy of the
# surrounding constant columns, with updated date / volume / prices).
$ws.Cells.Item(17, 1).Value = 4
$ws.Cells.Item(17, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(17, 3).Value = "Los Lagos"
$ws.Cells.Item(17, 4).Value = 44873
$ws.Cells.Item(17, 5).Value = 10
$ws.Cells.Item(17, 6).Value = "Fruta"
$ws.Cells.Item(17, 7).Value = 100107
$ws.Cells.Item(17, 8).Value = "Otros"
$ws.Cells.Item(17, 9).Value = 100107002
$ws.Cells.Item(17, 10).Value = "Chirimoya"
$ws.Cells.Item(17, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(17, 12).Value = "Primera"
$ws.Cells.Item(17, 13).Value = 300
$ws.Cells.Item(17, 14).Value = 22000
$ws.Cells.Item(17, 15).Value = 22500
$ws.Cells.Item(17, 16).Value = 22250
$ws.Cells.Item(17, 17).Value = "$/bandeja 8 kilos"
$ws.Cells.Item(17, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(17, 19).Value = 2781
$ws.Cells.Item(17, 20).Value = 8

# Make sure the date cell keeps the date-formatted style used by the
# other rows in column D.
$ws.Cells.Item(17, 4).NumberFormat = $ws.Cells.Item(18, 4).NumberFormat
